# Generate Report for Archive
# ------------------------------------------------------------------
# The nightly localization-status report flips each file's status from
# "Ready for handoff" to "In Translation" once the handoff package has
# actually been picked up by the translation vendor. Because the new
# status text is shorter, the Status columns (and their mirrored
# zh-cn/de-de columns on the Overview tab) are re-sized to fit the new
# content before the workbook is archived.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (E) / de-de (F) status columns for both rows ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Re-fit the now-narrower status columns to match the shorter text ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.56   # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = 12.56   # F: de-de status
$wsZhCn.Columns.Item(3).ColumnWidth = 12.56        # C: Status
$wsDeDe.Columns.Item(3).ColumnWidth = 12.56        # C: Status
